# Add new indicator-formula rows to the "Library_Formula" sheet.
# (This is the workbook's ActiveSheet / ActiveTab per xl/workbook.xml.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "INDICATOR_xx" formula names being registered, in the same order
# they are appended to the sheet (rows 100-148).
$indicators = @(
    71, 72, 73, 74, 75, 76, 77, 78, 79, 80,
    237, 238, 239, 240, 241, 242, 243, 244, 245, 246,
    247, 248, 249, 250, 251, 252, 253, 254, 255, 256,
    257, 258, 259, 260, 261, 262, 263, 264, 265, 266,
    267, 268, 269, 270, 271, 272, 273, 274, 275
)

$row = 100
foreach ($ind in $indicators) {
    # Insert a fresh row so it inherits the formatting of the row above
    # (style for column A, and for B/C/E/F) instead of the bare column
    # default style.
    $ws.Rows.Item($row).Insert()

    $ws.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($row, 2).Value = "LIB_EWS_IT"
    $ws.Cells.Item($row, 3).Value = "INDICATOR_$ind"
    $ws.Cells.Item($row, 5).Value = "String"
    $ws.Cells.Item($row, 6).Value = "String"

    $row++
}

# Update the view: scroll down to the newly added rows and move the
# active selection onto the new block.
$win = $excel.ActiveWindow
$win.ScrollRow = 128
$win.ScrollColumn = 1
[void]$ws.Range("C121").Select()
